$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7 to make room for the new item (CATAFLAM),
# shifting the existing data rows (FATROXIM, KETOLAC, NORHINOSE, سرنجات) down by one.
$ws.Rows("7:7").Insert()

# Copy the formatting of the row right below (now row 8, which holds the data that used
# to be in row 7) into the newly inserted row 7 so it matches the rest of the table.
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Recreate the merged cells for the new row 7 (matching the pattern used by every data row).
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

# Columns L, N and P store numeric-looking values as text - force text format so the
# values are not silently re-interpreted as numbers.
$ws.Range("L7").NumberFormat = "@"
$ws.Range("N7").NumberFormat = "@"
$ws.Range("P7").NumberFormat = "@"

# Fill in the new row 7 with the CATAFLAM item data.
$ws.Range("A7").Value2 = 1
$ws.Range("C7").Value2 = "CATAFLAM 75MG/3ML 6 AMP."
$ws.Range("H7").Value2 = "1:0"
$ws.Range("L7").Value2 = "1"
$ws.Range("N7").Value2 = "120.00"
$ws.Range("P7").Value2 = "19.2000"
$ws.Range("Q7").Value2 = "0:1"

# Renumber the "م" (index) column for the rows that shifted down.
$ws.Range("A8").Value2 = 2
$ws.Range("A9").Value2 = 3
$ws.Range("A10").Value2 = 4
$ws.Range("A11").Value2 = 5

# Update the سرنجات 3 سم row (now row 11): quantity sold / total price changed.
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value2 = "6.0000"
$ws.Range("Q11").Value2 = "3:0"

# Update the grand total (now on row 12, previously row 11).
$ws.Range("P12").Value2 = 258

# Update the generated timestamp in the footer (now row 13, previously row 12).
$ws.Range("A13").Value2 = "Friday, 15 August, 2025 3:04 PM"
